# Update "paises" (countries) COVID-19 data sheet and provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated timestamp title in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 17:16"

# Full refreshed country table (sorted by total cases desc), starting at row 4.
# Includes a newly reported country (Belice) appended as the final row (197).
$data = @"
China|81093|39|72703|5120|1749|9|3270
Italia|59138|0|7024|46638|3000|0|5476
Estados Unidos|39389|5843|187|38735|1040|48|467
España|33089|4321|3355|27528|2355|434|2206
Alemania|27546|2673|422|27009|23|21|115
Iran|23049|1411|8376|12861|0|127|1812
Francia|16689|671|2200|13815|1746|0|674
Corea del Sur|8961|64|3166|5684|59|7|111
Suiza|8547|1073|131|8298|141|20|118
Reino Unido|5837|154|135|5413|20|8|289
Paises Bajos|4749|545|2|4534|435|34|213
Austria|4075|493|9|4045|14|5|21
Belgica|3743|342|401|3254|322|13|88
Noruega|2547|162|6|2531|41|3|10
Portugal|2060|460|14|2023|47|9|23
Suecia|2046|112|16|2005|104|4|25
Australia|1717|108|88|1622|11|0|7
Brasil|1620|74|2|1593|18|0|25
Canada|1563|93|18|1524|1|1|21
Malasia|1518|212|159|1345|57|4|14
Dinamarca|1450|55|1|1425|55|11|24
Israel|1238|167|37|1200|24|0|1
Turquia|1236|0|0|1206|0|0|30
Chequia|1165|45|6|1158|19|0|1
Japon|1101|0|235|825|49|0|41
Ecuador|981|192|3|960|2|4|18
Irlanda|906|0|5|897|29|0|4
Pakistan|873|97|13|854|0|1|6
Luxemburgo|798|0|6|784|3|0|8
Chile|746|114|11|733|7|1|2
Tailandia|721|122|52|668|7|0|1
Crucero|712|0|567|137|15|0|8
Finlandia|700|74|10|689|13|0|1
Polonia|692|58|13|671|3|1|8
Grecia|624|0|19|588|18|2|17
Islandia|588|20|51|536|14|0|1
Indonesia|579|65|30|500|0|1|49
Rumania|576|143|73|498|15|2|5
Arabia Saudita|562|51|19|543|0|0|0
Singapur|509|54|152|355|14|0|2
Catar|494|0|33|461|6|0|0
India|467|71|34|424|0|2|9
Filipinas|462|82|18|411|1|8|33
Eslovenia|442|28|0|439|12|1|3
Rusia|438|71|17|420|0|0|1
Sudafrica|402|128|4|398|0|0|0
Peru|363|0|1|357|5|0|5
Hong Kong|356|38|100|252|4|0|4
Estonia|352|26|4|348|4|0|0
Barein|339|5|160|177|3|0|2
Egipto|327|0|56|257|0|0|14
Mexico|316|65|4|310|1|0|2
Croacia|315|61|5|309|5|0|1
Panama|313|0|1|309|7|0|3
Libano|267|19|8|255|4|0|4
Argentina|266|0|27|235|0|0|4
Irak|266|33|62|181|0|3|23
Serbia|249|27|3|244|12|0|2
Republica Dominicana|245|43|0|242|0|0|3
Colombia|235|4|3|229|0|1|3
Argelia|230|29|65|148|0|0|17
Bulgaria|201|14|3|195|8|0|3
Emiratos Arabes Unidos|198|45|41|155|2|0|2
Taiwan|195|26|28|165|0|0|2
Armenia|194|0|2|192|6|0|0
Kuwait|189|1|30|159|5|0|0
San Marino|187|12|4|163|13|0|20
Eslovaquia|186|1|7|179|2|0|0
Letonia|180|41|1|179|0|0|0
Hungria|167|36|16|144|6|1|7
Lituania|160|17|1|158|1|0|1
Uruguay|158|0|0|158|2|0|0
Republica de Macedonia|136|21|1|133|1|1|2
Costa Rica|134|0|2|130|2|0|2
Marruecos|134|19|3|127|1|0|4
Principado de Andorra|133|20|1|131|2|0|1
Bosnia y Herzegovina|130|4|2|127|1|0|1
Vietnam|123|10|17|106|2|0|0
Islas Feroe|118|3|14|104|0|0|0
Jordania|112|0|1|111|0|0|0
Moldavia|109|15|2|106|3|0|1
Malta|107|17|2|105|1|0|0
Albania|104|15|2|98|2|2|4
Nueva Zelanda|102|0|0|102|0|0|0
Burkina Faso|99|24|5|90|0|0|4
Republica de Chipre|95|0|3|91|3|0|1
Sri Lanka|94|12|3|91|2|0|0
Brunei|91|3|2|89|2|0|0
Tunez|89|14|1|85|11|0|3
Camboya|86|2|2|84|1|0|0
Bielorrusia|81|5|22|59|0|0|0
Venezuela|77|7|15|62|2|0|0
Ucrania|73|0|1|69|0|0|3
Azerbaiyan|72|7|10|61|3|0|1
Reunion|71|7|1|70|0|0|0
Senegal|67|0|5|62|0|0|0
Oman|66|11|17|49|0|0|0
Kazajistan|62|2|0|62|0|0|0
Georgia|61|7|8|53|1|0|0
Estado de Palestina|59|0|17|42|0|0|0
Guadalupe|58|0|0|57|4|0|1
Camerun|56|16|2|54|0|0|0
Trinidad yTobago|51|1|0|51|0|0|0
Uzbekistan|46|3|0|46|0|0|0
Liechtenstein|46|9|0|46|0|0|0
Martinica|44|0|0|43|7|0|1
Afganistan|40|0|1|38|0|0|1
Consejo Danes para los Refugiados|36|6|0|35|0|0|1
Nigeria|36|6|2|33|0|1|1
Cuba|35|0|0|34|0|0|1
Banglades|33|6|5|25|0|1|3
Puerto Rico|31|8|1|28|0|1|2
Guam|29|2|0|28|0|0|1
Mauricio|28|0|0|26|1|0|2
Bolivia|27|3|0|27|0|0|0
Honduras|27|1|0|27|0|0|0
Costa de Marfil|25|11|2|23|0|0|0
Ghana|24|1|0|23|0|0|1
Macao|24|2|10|14|0|0|0
Monaco|23|0|1|22|0|0|0
Montenegro|22|1|0|21|0|1|1
Paraguay|22|0|0|21|1|0|1
Mayotte|21|10|0|21|0|0|0
Guayana Francesa|20|2|6|14|0|0|0
Ruanda|19|0|0|19|0|0|0
Guyana|19|0|0|18|0|0|1
Guatemala|19|0|0|18|0|0|1
Jamaica|19|0|2|16|0|0|1
Polinesia Francesa|18|0|0|18|0|0|0
Togo|18|2|0|18|0|0|0
Barbados|17|3|0|17|0|0|0
Islas Virgenes de los Estados Unidos|17|11|0|17|0|0|0
Kenia|16|1|0|16|0|0|0
Kirguistan|16|2|0|16|0|0|0
Gibraltar|15|0|5|10|0|0|0
Isla de Man|13|8|0|13|0|0|0
Maldivas|13|0|5|8|0|0|0
Tanzania|12|0|0|12|0|0|0
Madagascar|12|9|0|12|0|0|0
Etiopia|11|0|0|11|0|0|0
Mongolia|10|0|0|10|0|0|0
Guinea Ecuatorial|9|3|0|9|0|0|0
Aruba|9|0|1|8|0|0|0
Nueva Caledonia|8|4|0|8|0|0|0
Seychelles|7|0|0|7|0|0|0
Bermudas|6|0|0|6|0|0|0
San Martin (Parte Francesa)|5|0|0|5|0|0|0
Surinam|5|0|0|5|0|0|0
Haiti|5|3|0|5|0|0|0
Gabon|5|0|0|4|0|0|1
Groenlandia|4|0|0|4|0|0|0
Suazilandia|4|0|0|4|0|0|0
Bahamas|4|0|0|4|0|0|0
Guinea|4|2|0|4|0|0|0
Curazao|4|1|0|3|0|0|1
Congo|3|0|0|3|0|0|0
Zambia|3|0|0|3|0|0|0
Cabo Verde|3|0|0|3|0|0|0
Fiyi|3|1|0|3|0|0|0
Namibia|3|0|0|3|0|0|0
San Bartolome|3|0|0|3|0|0|0
Liberia|3|0|0|3|0|0|0
El Salvador|3|0|0|3|0|0|0
Republica de Africa Central|3|0|0|3|0|0|0
Islas Caimanes|3|0|0|2|0|0|1
Zimbabue|3|0|0|2|0|1|1
Santa Lucia|2|0|0|2|0|0|0
Butan|2|0|0|2|0|0|0
Benin|2|0|0|2|0|0|0
Mauritania|2|0|0|2|0|0|0
Angola|2|0|0|2|0|0|0
Nicaragua|2|0|0|2|0|0|0
Niger|2|0|0|2|0|0|0
Nepal|2|1|1|1|0|0|0
Gambia|2|1|0|1|0|1|1
Sudan|2|0|0|1|0|0|1
Mozambique|1|0|0|1|0|0|0
Siria|1|0|0|1|0|0|0
Republica del Chad|1|0|0|1|0|0|0
San Martin (Parte Holandesa)|1|0|0|1|0|0|0
Republica de Yibuti|1|0|0|1|0|0|0
Granada|1|0|0|1|0|0|0
San Vicente y las Granadinas|1|0|0|1|0|0|0
Dominica|1|0|0|1|0|0|0
Islas Turcas y Caicos|1|1|0|1|0|0|0
Montserrat|1|0|0|1|0|0|0
Somalia|1|0|0|1|0|0|0
Santa Sede|1|0|0|1|0|0|0
Belice|1|1|0|1|0|0|0
Papua Nueva Guinea|1|0|0|1|0|0|0
Eritrea|1|0|0|1|0|0|0
Timor Oriental|1|0|0|1|0|0|0
Antigua y Barbuda|1|0|0|1|0|0|0
Uganda|1|0|0|1|0|0|0
"@

$lines = $data -split "`n"
$r = 4
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $ws.Cells.Item($r, 1).Value = $parts[0]
    $ws.Cells.Item($r, 2).Value = [double]$parts[1]
    $ws.Cells.Item($r, 3).Value = [double]$parts[2]
    $ws.Cells.Item($r, 4).Value = [double]$parts[3]
    $ws.Cells.Item($r, 5).Value = [double]$parts[4]
    $ws.Cells.Item($r, 6).Value = [double]$parts[5]
    $ws.Cells.Item($r, 7).Value = [double]$parts[6]
    $ws.Cells.Item($r, 8).Value = [double]$parts[7]
    $r++
}

Write-Output "Updated rows 4..$($r-1)"
